$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. Insert a fresh row at 86,
# pushing the existing rows 86-111 down to 87-112, then populate the new
# row with the new data point (dimension grows from A1:T111 to A1:T112).
$ws.Rows.Item(86).Insert()

$ws.Cells.Item(86, 1).Value = 10
$ws.Cells.Item(86, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(86, 3).Value = "La Araucanía"
$ws.Cells.Item(86, 4).Value = 44900
$ws.Cells.Item(86, 5).Value = 9
$ws.Cells.Item(86, 6).Value = "Fruta"
$ws.Cells.Item(86, 7).Value = 100101
$ws.Cells.Item(86, 8).Value = "Berries"
$ws.Cells.Item(86, 9).Value = 100101001
$ws.Cells.Item(86, 10).Value = "Arándano (blue)"
$ws.Cells.Item(86, 11).Value = "Sin especificar"
$ws.Cells.Item(86, 12).Value = "Primera"
$ws.Cells.Item(86, 13).Value = 1100
$ws.Cells.Item(86, 14).Value = 1800
$ws.Cells.Item(86, 15).Value = 2000
$ws.Cells.Item(86, 16).Value = 1891
$ws.Cells.Item(86, 17).Value = "`$/kilo"
$ws.Cells.Item(86, 18).Value = "Región del Maule"
$ws.Cells.Item(86, 19).Value = 1891
$ws.Cells.Item(86, 20).Value = 1
